# Fix Training Data Issue: the "Date" column (BF) held a malformed literal
# like "6-1-2007-08" instead of the actual game date "2008-06-01".
# Data was taken from 1 day off due to way NBA stats were shown.
#
# The corrected value must stay literal text (it is not a real Excel date),
# so it is entered with a leading apostrophe - exactly like a user typing
# '2008-06-01 into the cell - which keeps it as text instead of having
# Excel reinterpret the ISO-looking string as a date serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $ws.Cells.Item($row, 58).Value = "'2008-06-01"
}
